$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired order of player rows (row 1 is the header and is unchanged).
# The edit re-sorts the roster: "Bradley Beal" moves from row 4 down to the
# bottom of the list (row 18); every other player shifts up to fill the gap.
$data = @(
    @("Jose Alvarado", "PG", "New Orleans Pelicans"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Donte DiVincenzo", "PG,SG,SF", "Minnesota Timberwolves"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Draymond Green", "PF,C", "Golden State Warriors"),
    @("Amen Thompson", "SG,SF,PF", "Houston Rockets"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
